# Insert a new weekly record at row 188 of the daily-price sheet.
# This shifts the existing rows 188-298 down to 189-299 (dimension
# grows from A1:R298 to A1:R299) and populates the new row with the
# latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 188..298 down by one to make room for the new record.
$ws.Rows(188).Insert()

# Populate the newly inserted row 188 with the new weekly data point.
$ws.Range("A188").Value = 3
$ws.Range("B188").Value = 'Femacal de La Calera'
$ws.Range("C188").Value = 'Coquimbo'
$ws.Range("D188").Value = 44719
$ws.Range("E188").Value = 5
$ws.Range("F188").Value = 100112001
$ws.Range("G188").Value = 'Berenjena'
$ws.Range("H188").Value = 'Sin especificar'
$ws.Range("I188").Value = 'Primera'
$ws.Range("J188").Value = 105
$ws.Range("K188").Value = 6000
$ws.Range("L188").Value = 7000
$ws.Range("M188").Value = 6476
$ws.Range("N188").Value = '$/caja 60 unidades'
$ws.Range("O188").Value = 'Región de Arica y Parinacota'
$ws.Range("P188").Value = 108
$ws.Range("Q188").Value = 60
$ws.Range("R188").Value = 'Hortaliza'
